# ShapeCrawler autoshape-case003.pptx update:
#  1. Re-cache the "datetimeFigureOut" date placeholder text (master + every
#     layout) from 11.10.2022 to 01.11.2022 - this happens automatically in
#     PowerPoint whenever the deck is touched/resaved on a different day.
#  2. Add a third AutoShape ("AutoShape 3") to slide 1, to the right of the
#     existing "AutoShape 2", containing the text "Some text" (two runs with
#     different font sizes), used as a fixture for
#     IParagraph.ReplaceText(string oldValue, string newValue).

$p = $ppt.ActivePresentation

# --- 1. Update every "Date Placeholder *" shape's cached field text ---
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "01.11.2022"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholders $layout.Shapes
}

# --- 2. Add the new "AutoShape 3" shape to slide 1 ---
$slide = $p.Slides.Item(1)

# Burn a throwaway shape id first so the real shape lands on id=5, matching
# how id=3 is already missing/consumed in the source deck (AutoShape 2 = 2,
# AutoShape 1 = 4).
$placeholderShape = $slide.Shapes.AddShape(1, 0, 0, 1, 1)
$placeholderShape.Delete()

# Clone "AutoShape 2" (id=2) - same size, noFill + accent1 shape style -
# then reposition/rename/retext it into "AutoShape 3".
$newShape = $slide.Shapes.Item(1).Duplicate()
$newShape.Name = "AutoShape 3"
$newShape.Left = 2837089 / 12700
$newShape.Top = 159201 / 12700

$tr = $newShape.TextFrame.TextRange
$tr.Text = "Some text"
$tr.Characters(1, 4).Font.Size = 12
$tr.Characters(5, 5).Font.Size = 11
